$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert the three new rows (shifts everything below down)
$ws.Range("A4").EntireRow.Insert()
$ws.Range("A9").EntireRow.Insert()
$ws.Range("A10").EntireRow.Insert()

# New "ngx-toastr" entry (row 9)
$ws.Hyperlinks.Add($ws.Range("C9"), "https://www.npmjs.com/package/ngx-toastr") | Out-Null
$ws.Range("B9").Value = "ngx-toastr"

# New "auth0/angular-jwt" entry (row 10)
$ws.Range("B10").Value = "auth0/angular-jwt"
$ws.Range("C10").Value = "This library provides an HttpInterceptor which automatically attaches a JSON Web Token to HttpClient requests."
$ws.Range("D10").Value = "https://www.npmjs.com/package/@auth0/angular-jwt"

# New "Bootswatch" entry (row 4)
$ws.Hyperlinks.Add($ws.Range("C4"), "https://bootswatch.com/") | Out-Null
$ws.Range("D4").Value = "Free themes for Bootstrap"
$ws.Range("B4").Value = "Bootswatch"

# Reset the view: scroll back to top-left and select C2
$ws.Range("C2").Select()

Write-Host "done"
